$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D as text so numeric-looking price strings (e.g. "1.002")
# are not auto-converted to numbers by Excel's type inference on assignment.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "24.156.58"
$ws.Range("D3").Value = "1.644.66"
$ws.Range("E3").Value = "  -3.27%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "308.11"
$ws.Range("E5").Value = "  -2.40%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "0.3913"
$ws.Range("E7").Value = "  -1.63%  "
$ws.Range("D8").Value = "0.3862"
$ws.Range("E8").Value = "  -4.02%  "
$ws.Range("D9").Value = "1.002"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").Value = "1.358"
$ws.Range("E10").Value = "  -7.64%  "
$ws.Range("D11").Value = "48.59"
$ws.Range("E11").Value = "  -8.71%  "
$ws.Range("D12").Value = "0.08471"
$ws.Range("E12").Value = "  -3.82%  "
$ws.Range("D13").Value = "24.09"
$ws.Range("E13").Value = "  -7.20%  "
$ws.Range("D14").Value = "7.158"
$ws.Range("E14").Value = "  -4.11%  "
$ws.Range("D15").Value = "0.00001286"
$ws.Range("E15").Value = "  -4.90%  "
$ws.Range("D16").Value = "7.504"
$ws.Range("E16").Value = "  -5.68%  "
$ws.Range("D17").Value = "1.645.94"
$ws.Range("E17").Value = "  -3.52%  "
$ws.Range("D18").Value = "94.40"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("D19").Value = "0.06954"
$ws.Range("E19").Value = "  -3.26%  "
$ws.Range("D20").Value = "20.98"
$ws.Range("E20").Value = "  +1.52%  "
$ws.Range("D21").Value = "6.966"
$ws.Range("E21").Value = "  -4.82%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").Value = "13.77"
$ws.Range("E23").Value = "  -4.11%  "
$ws.Range("D24").Value = "24.147.86"
$ws.Range("E24").Value = "  -3.30%  "
$ws.Range("E25").Value = "  -1.52%  "
$ws.Range("D26").Value = "2.734"
$ws.Range("E26").Value = "  -7.29%  "
$ws.Range("D27").Value = "22.54"
$ws.Range("E27").Value = "  -4.88%  "
$ws.Range("D28").Value = "8.940"
$ws.Range("E28").Value = "  +6.94%  "
$ws.Range("D29").Value = "157.91"
$ws.Range("E29").Value = "  -2.95%  "
$ws.Range("D30").Value = "141.75"
$ws.Range("E30").Value = "  -6.10%  "
$ws.Range("D31").Value = "5.367"
$ws.Range("E31").Value = "  -13.22%  "
$ws.Range("D32").Value = "2.471"
$ws.Range("E32").Value = "  -6.51%  "
$ws.Range("D33").Value = "1.828.34"
$ws.Range("E33").Value = "  -3.91%  "
$ws.Range("D34").Value = "7.259"
$ws.Range("E34").Value = "  +1.42%  "
$ws.Range("D35").Value = "0.08056"
$ws.Range("E35").Value = "  -5.72%  "
$ws.Range("D36").Value = "0.9823"
$ws.Range("E36").Value = "  -5.33%  "
$ws.Range("D37").Value = "0.02941"
$ws.Range("E37").Value = "  -6.69%  "
$ws.Range("D38").Value = "0.2715"
$ws.Range("E38").Value = "  -5.30%  "
$ws.Range("D39").Value = "0.09250"
$ws.Range("E39").Value = "  -3.45%  "
$ws.Range("D40").Value = "1.479"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("D41").Value = "10.11"
$ws.Range("E41").Value = "  -7.13%  "
$ws.Range("D42").Value = "0.7651"
$ws.Range("E42").Value = "  -7.38%  "
$ws.Range("D43").Value = "13.15"
$ws.Range("E43").Value = "  -6.20%  "
$ws.Range("D44").Value = "16.04"
$ws.Range("E44").Value = "  -6.41%  "
$ws.Range("D45").Value = "2.492"
$ws.Range("E45").Value = "  -7.31%  "
$ws.Range("D46").Value = "0.6896"
$ws.Range("E46").Value = "  -6.72%  "
$ws.Range("D47").Value = "4.096"
$ws.Range("E47").Value = "  -3.47%  "
$ws.Range("D48").Value = "1.002"
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("D49").Value = "0.08413"
$ws.Range("E49").Value = "  -4.39%  "
$ws.Range("D50").Value = "134.32"
$ws.Range("E50").Value = "  -3.52%  "
$ws.Range("D51").Value = "1.265"
$ws.Range("E51").Value = "  -9.15%  "

# Remove the temporary text-number-format so styling matches the original
# (values remain text; only the now-unneeded explicit format is cleared).
$ws.Range("D2:D51").ClearFormats()
